$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-05-15 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-16 Thursday", 2) | Out-Null

# Update each answer cell in the single 20-row x 5-col table.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Find.Execute("51÷7=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=2, 2", 2) | Out-Null
$t.Cell(1, 2).Range.Find.Execute("83÷6=13, 5", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=1, 7", 2) | Out-Null
$t.Cell(1, 3).Range.Find.Execute("67÷6=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "10÷2=5, 0", 2) | Out-Null
$t.Cell(1, 4).Range.Find.Execute("80÷5=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "92÷2=46, 0", 2) | Out-Null
$t.Cell(1, 5).Range.Find.Execute("56÷8=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=32, 2", 2) | Out-Null
$t.Cell(5, 1).Range.Find.Execute("54÷3=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "30÷2=15, 0", 2) | Out-Null
$t.Cell(5, 2).Range.Find.Execute("67÷5=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "33÷2=16, 1", 2) | Out-Null
$t.Cell(5, 3).Range.Find.Execute("32÷8=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "39÷6=6, 3", 2) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("98÷8=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷9=3, 8", 2) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("59÷5=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "36÷3=12, 0", 2) | Out-Null
$t.Cell(9, 1).Range.Find.Execute("36÷5=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "51÷5=10, 1", 2) | Out-Null
$t.Cell(9, 2).Range.Find.Execute("59÷6=9, 5", $true, $false, $false, $false, $false, $true, 1, $false, "22÷4=5, 2", 2) | Out-Null
$t.Cell(9, 3).Range.Find.Execute("43÷5=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "91÷7=13, 0", 2) | Out-Null
$t.Cell(9, 4).Range.Find.Execute("27÷6=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "76÷3=25, 1", 2) | Out-Null
$t.Cell(9, 5).Range.Find.Execute("95÷8=11, 7", $true, $false, $false, $false, $false, $true, 1, $false, "77÷4=19, 1", 2) | Out-Null
$t.Cell(13, 1).Range.Find.Execute("53÷8=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "77÷6=12, 5", 2) | Out-Null
$t.Cell(13, 2).Range.Find.Execute("60÷8=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=39, 1", 2) | Out-Null
$t.Cell(13, 3).Range.Find.Execute("25÷4=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "69÷4=17, 1", 2) | Out-Null
$t.Cell(13, 4).Range.Find.Execute("13÷5=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷3=13, 0", 2) | Out-Null
$t.Cell(13, 5).Range.Find.Execute("58÷4=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "72÷7=10, 2", 2) | Out-Null
$t.Cell(17, 1).Range.Find.Execute("46÷7=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "81÷4=20, 1", 2) | Out-Null
$t.Cell(17, 2).Range.Find.Execute("87÷8=10, 7", $true, $false, $false, $false, $false, $true, 1, $false, "76÷3=25, 1", 2) | Out-Null
$t.Cell(17, 3).Range.Find.Execute("17÷5=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "37÷3=12, 1", 2) | Out-Null
$t.Cell(17, 4).Range.Find.Execute("81÷7=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "74÷8=9, 2", 2) | Out-Null
$t.Cell(17, 5).Range.Find.Execute("71÷7=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "80÷5=16, 0", 2) | Out-Null
